# Convert the two M2Doc "documentKeywords" field codes (fldChar begin/instrText.../fldChar end)
# into plain literal text runs "{...}" (TokenIteratorFieldRewriterSplit no longer needs real
# Word fields - it rewrites the template text directly).
#
# For each target paragraph we replace the run content (fldChar begin .. fldChar end) with a
# fresh set of <w:r><w:t>...</w:t></w:r> runs carrying the same characters, except the field's
# leading space becomes "{" and its trailing space becomes "}". The paragraph mark itself (and
# its <w:pPr>) is left untouched by stopping the replaced Range one character before the
# paragraph end, and the original <w:p> identity attributes are reasserted explicitly in the
# injected OOXML so they survive the InsertXML round-trip.

$d = $word.ActiveDocument

function Convert-FieldParagraphToText($Paragraph, $ParaId, $TextId, $RsidR, $RsidP, $Texts) {
    $runsXml = ""
    foreach ($t in $Texts) {
        $escaped = $t.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $runsXml += "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    }

    $xml = "<?xml version=`"1.0`" standalone=`"yes`"?>" +
        "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
        "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
        "<pkg:xmlData>" +
        "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" xmlns:w14=`"http://schemas.microsoft.com/office/word/2010/wordml`" xmlns:mc=`"http://schemas.openxmlformats.org/markup-compatibility/2006`" mc:Ignorable=`"w14`">" +
        "<w:body>" +
        "<w:p w14:paraId=`"$ParaId`" w14:textId=`"$TextId`" w:rsidR=`"$RsidR`" w:rsidRDefault=`"$RsidR`" w:rsidP=`"$RsidP`">" +
        "<w:pPr><w:tabs><w:tab w:val=`"left`" w:pos=`"3119`"/></w:tabs><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
        $runsXml +
        "</w:p>" +
        "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

    $full = $Paragraph.Range
    $target = $d.Range($full.Start, $full.End - 1)
    $target.InsertXML($xml)
}

# Paragraph 2: { m:'Some value'.setDocumentKeywords() }  ->  {m:'Some value'.setDocumentKeywords()}
Convert-FieldParagraphToText $d.Paragraphs.Item(2) "5F0A223D" "7BAB5916" "00E1471F" "00E1471F" @("{m:", "'", "Some value", "'", ".", "setDocument", "Keywords", "()}")

# Paragraph 3: { m:''.getDocumentKeywords() }  ->  {m:''.getDocumentKeywords()}
Convert-FieldParagraphToText $d.Paragraphs.Item(3) "2C980985" "50C6D369" "00CD75A1" "00CD75A1" @("{m:''.g", "etDocument", "Keywor", "d", "s", "()}")
